$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 1.551830800456742
$ws.Cells.Item(2, 3).Value = 0.1837026811945179
$ws.Cells.Item(2, 5).Value = 0.08309285813966838
$ws.Cells.Item(2, 6).Value = 0.4443680307746263
$ws.Cells.Item(2, 7).Value = 0.9951647432747848
$ws.Cells.Item(2, 8).Value = 0.9768654957502179
$ws.Cells.Item(2, 12).Value = 0.2275577897234911
$ws.Cells.Item(2, 14).Value = 1.303212590595503
$ws.Cells.Item(3, 2).Value = 1.435940737329986
$ws.Cells.Item(3, 3).Value = 0.16798027384948
$ws.Cells.Item(3, 5).Value = 0.0834689652082119
$ws.Cells.Item(3, 6).Value = 0.387822817061874
$ws.Cells.Item(3, 7).Value = 0.9820497032226712
$ws.Cells.Item(3, 8).Value = 0.9771427098800416
$ws.Cells.Item(3, 12).Value = 0.218211674527879
$ws.Cells.Item(3, 14).Value = 1.322286152871644
$ws.Cells.Item(4, 2).Value = 1.36544419479003
$ws.Cells.Item(4, 3).Value = 0.1582513534000896
$ws.Cells.Item(4, 5).Value = 0.08373409204934035
$ws.Cells.Item(4, 6).Value = 0.3531389305168915
$ws.Cells.Item(4, 7).Value = 0.9748565468470929
$ws.Cells.Item(4, 8).Value = 0.9779422425865079
$ws.Cells.Item(4, 12).Value = 0.2126037950455952
$ws.Cells.Item(4, 14).Value = 1.334597215390133
$ws.Cells.Item(5, 2).Value = 1.336882217070411
$ws.Cells.Item(5, 3).Value = 0.1542677469071521
$ws.Cells.Item(5, 5).Value = 0.0838507359652283
$ws.Cells.Item(5, 6).Value = 0.3390132514313251
$ws.Cells.Item(5, 7).Value = 0.9721402288361958
$ws.Cells.Item(5, 8).Value = 0.9784258916045587
$ws.Cells.Item(5, 12).Value = 0.2103513370967249
$ws.Cells.Item(5, 14).Value = 1.339764741325979
$ws.Cells.Item(6, 2).Value = 1.332149541195975
$ws.Cells.Item(6, 3).Value = 0.1536051233490809
$ws.Cells.Item(6, 5).Value = 0.08387062428842817
$ws.Cells.Item(6, 6).Value = 0.336668177824194
$ws.Cells.Item(6, 7).Value = 0.9717021343146257
$ws.Cells.Item(6, 8).Value = 0.9785157218053797
$ws.Cells.Item(6, 12).Value = 0.2099792974045016
$ws.Cells.Item(6, 14).Value = 1.340631900430679
$ws.Cells.Item(7, 2).Value = 1.365058326268468
$ws.Cells.Item(7, 3).Value = 0.1581977060626798
$ws.Cells.Item(7, 5).Value = 0.08373563031314823
$ws.Cells.Item(7, 6).Value = 0.3529483938344953
$ws.Cells.Item(7, 7).Value = 0.9748190448719356
$ws.Cells.Item(7, 8).Value = 0.9779481267435131
$ws.Cells.Item(7, 12).Value = 0.2125732848872133
$ws.Cells.Item(7, 14).Value = 1.334666296622901
$ws.Cells.Item(8, 2).Value = 1.511734711314375
$ws.Cells.Item(8, 3).Value = 0.178297183753358
$ws.Cells.Item(8, 5).Value = 0.08321544474772757
$ws.Cells.Item(8, 6).Value = 0.4248636149813336
$ws.Cells.Item(8, 7).Value = 0.9904635604340939
$ws.Cells.Item(8, 8).Value = 0.9768302097231185
$ws.Cells.Item(8, 12).Value = 0.2243080877693586
$ws.Cells.Item(8, 14).Value = 1.309664566548934
$ws.Cells.Item(9, 2).Value = 1.804629888338013
$ws.Cells.Item(9, 3).Value = 0.2171191792831166
$ws.Cells.Item(9, 5).Value = 0.08246656499196092
$ws.Cells.Item(9, 6).Value = 0.5661985755041457
$ws.Cells.Item(9, 7).Value = 1.028019046609387
$ws.Cells.Item(9, 8).Value = 0.9796515240484212
$ws.Cells.Item(9, 12).Value = 0.2483610955125215
$ws.Cells.Item(9, 14).Value = 1.265402125832015
$ws.Cells.Item(10, 2).Value = 2.023085094830378
$ws.Cells.Item(10, 3).Value = 0.2452889679385919
$ws.Cells.Item(10, 5).Value = 0.08208159403843851
$ws.Cells.Item(10, 6).Value = 0.6702781546542269
$ws.Cells.Item(10, 7).Value = 1.059886297703201
$ws.Cells.Item(10, 8).Value = 0.9848107095034493
$ws.Cells.Item(10, 12).Value = 0.2666756477086665
$ws.Cells.Item(10, 14).Value = 1.235796059257627
$ws.Cells.Item(11, 2).Value = 2.123188961850303
$ws.Cells.Item(11, 3).Value = 0.2580296139829841
$ws.Cells.Item(11, 5).Value = 0.08194232541086066
$ws.Cells.Item(11, 6).Value = 0.7176906081379002
$ws.Cells.Item(11, 7).Value = 1.075330142148317
$ws.Cells.Item(11, 8).Value = 0.9878345973293108
$ws.Cells.Item(11, 12).Value = 0.2751489058414478
$ws.Cells.Item(11, 14).Value = 1.222961972170143
$ws.Cells.Item(12, 2).Value = 2.161200941578443
$ws.Cells.Item(12, 3).Value = 0.2628436449359413
$ws.Cells.Item(12, 5).Value = 0.08189474254502649
$ws.Cells.Item(12, 6).Value = 0.7356546913071611
$ws.Cells.Item(12, 7).Value = 1.081315946291852
$ws.Cells.Item(12, 8).Value = 0.9890775251129185
$ws.Cells.Item(12, 12).Value = 0.2783780216711023
$ws.Cells.Item(12, 14).Value = 1.218193387922286
$ws.Cells.Item(13, 2).Value = 2.153009716341501
$ws.Cells.Item(13, 3).Value = 0.2618073266555712
$ws.Cells.Item(13, 5).Value = 0.08190476110440414
$ws.Cells.Item(13, 6).Value = 0.7317853510981394
$ws.Cells.Item(13, 7).Value = 1.080020657034538
$ws.Cells.Item(13, 8).Value = 0.9888054785345446
$ws.Cells.Item(13, 12).Value = 0.2776816618674189
$ws.Cells.Item(13, 14).Value = 1.219216318375672
$ws.Cells.Item(14, 2).Value = 2.12631412834935
$ws.Cells.Item(14, 3).Value = 0.2584258788204181
$ws.Cells.Item(14, 5).Value = 0.08193830742822072
$ws.Cells.Item(14, 6).Value = 0.7191683204515869
$ws.Cells.Item(14, 7).Value = 1.075819833340205
$ws.Cells.Item(14, 8).Value = 0.9879348898461444
$ws.Cells.Item(14, 12).Value = 0.2754141563189449
$ws.Cells.Item(14, 14).Value = 1.222567824462116
$ws.Cells.Item(15, 2).Value = 2.109975966684317
$ws.Cells.Item(15, 3).Value = 0.2563532688826626
$ws.Cells.Item(15, 5).Value = 0.08195952684488184
$ws.Cells.Item(15, 6).Value = 0.7114413442032514
$ws.Cells.Item(15, 7).Value = 1.073264661202643
$ws.Cells.Item(15, 8).Value = 0.9874143867282896
$ws.Cells.Item(15, 12).Value = 0.2740279126578855
$ws.Cells.Item(15, 14).Value = 1.224632630234524
$ws.Cells.Item(16, 2).Value = 2.016557727825898
$ws.Cells.Item(16, 3).Value = 0.2444548538085201
$ws.Cells.Item(16, 5).Value = 0.08209141700850076
$ws.Cells.Item(16, 6).Value = 0.6671810134426437
$ws.Cells.Item(16, 7).Value = 1.058896181921568
$ws.Cells.Item(16, 8).Value = 0.9846267561070476
$ws.Cells.Item(16, 12).Value = 0.266124761661942
$ws.Cells.Item(16, 14).Value = 1.23664757379033
$ws.Cells.Item(17, 2).Value = 1.959435132744943
$ws.Cells.Item(17, 3).Value = 0.2371366650866378
$ws.Cells.Item(17, 5).Value = 0.0821815103274588
$ws.Cells.Item(17, 6).Value = 0.6400460337215605
$ws.Cells.Item(17, 7).Value = 1.050325154064694
$ws.Cells.Item(17, 8).Value = 0.9830903737374399
$ws.Cells.Item(17, 12).Value = 0.2613128253372707
$ws.Cells.Item(17, 14).Value = 1.244180880510525
$ws.Cells.Item(18, 2).Value = 1.926648188924219
$ws.Cells.Item(18, 3).Value = 0.2329204741669457
$ws.Cells.Item(18, 5).Value = 0.08223670466388455
$ws.Cells.Item(18, 6).Value = 0.6244449056556647
$ws.Cells.Item(18, 7).Value = 1.045484406683727
$ws.Cells.Item(18, 8).Value = 0.9822703744617343
$ws.Cells.Item(18, 12).Value = 0.2585584762981625
$ws.Cells.Item(18, 14).Value = 1.248573485469926
$ws.Cells.Item(19, 2).Value = 1.915558845629334
$ws.Cells.Item(19, 3).Value = 0.2314917469044815
$ws.Cells.Item(19, 5).Value = 0.08225597221730574
$ws.Cells.Item(19, 6).Value = 0.6191636801734006
$ws.Cells.Item(19, 7).Value = 1.04386067191632
$ws.Cells.Item(19, 8).Value = 0.982003659900414
$ws.Cells.Item(19, 12).Value = 0.257628191565118
$ws.Cells.Item(19, 14).Value = 1.250070982581914
$ws.Cells.Item(20, 2).Value = 1.965508843534451
$ws.Cells.Item(20, 3).Value = 0.2379164189367486
$ws.Cells.Item(20, 5).Value = 0.08217157044376755
$ws.Cells.Item(20, 6).Value = 0.642933953830422
$ws.Cells.Item(20, 7).Value = 1.051228325400672
$ws.Cells.Item(20, 8).Value = 0.983247329275315
$ws.Cells.Item(20, 12).Value = 0.2618236820144233
$ws.Cells.Item(20, 14).Value = 1.243372772177848
$ws.Cells.Item(21, 2).Value = 2.134152423194337
$ws.Cells.Item(21, 3).Value = 0.2594193785669745
$ws.Cells.Item(21, 5).Value = 0.08192831415786017
$ws.Cells.Item(21, 6).Value = 0.7228739723492197
$ws.Cells.Item(21, 7).Value = 1.077049973156988
$ws.Cells.Item(21, 8).Value = 0.9881879432716119
$ws.Cells.Item(21, 12).Value = 0.2760796213176491
$ws.Cells.Item(21, 14).Value = 1.221580923061953
$ws.Cells.Item(22, 2).Value = 2.244981819729787
$ws.Cells.Item(22, 3).Value = 0.2734111929060816
$ws.Cells.Item(22, 5).Value = 0.08179937987323527
$ws.Cells.Item(22, 6).Value = 0.7751780083420101
$ws.Cells.Item(22, 7).Value = 1.094728319012319
$ws.Cells.Item(22, 8).Value = 0.9919874850419887
$ws.Cells.Item(22, 12).Value = 0.2855160930084395
$ws.Cells.Item(22, 14).Value = 1.207871598563752
$ws.Cells.Item(23, 2).Value = 2.185774099210278
$ws.Cells.Item(23, 3).Value = 0.2659491167671888
$ws.Cells.Item(23, 5).Value = 0.08186544544132701
$ws.Cells.Item(23, 6).Value = 0.7472568307915566
$ws.Cells.Item(23, 7).Value = 1.085219188058403
$ws.Cells.Item(23, 8).Value = 0.9899072238304711
$ws.Cells.Item(23, 12).Value = 0.2804687207756587
$ws.Cells.Item(23, 14).Value = 1.215139674935994
$ws.Cells.Item(24, 2).Value = 1.962762751956461
$ws.Cells.Item(24, 3).Value = 0.2375639198463944
$ws.Cells.Item(24, 5).Value = 0.08217605367947023
$ws.Cells.Item(24, 6).Value = 0.6416283278902313
$ws.Cells.Item(24, 7).Value = 1.050819731211817
$ws.Cells.Item(24, 8).Value = 0.9831761725773447
$ws.Cells.Item(24, 12).Value = 0.2615926860371474
$ws.Cells.Item(24, 14).Value = 1.243737926071777
$ws.Cells.Item(25, 2).Value = 1.72482475039601
$ws.Cells.Item(25, 3).Value = 0.2066793862462077
$ws.Cells.Item(25, 5).Value = 0.0826401337680096
$ws.Cells.Item(25, 6).Value = 0.5279251897347308
$ws.Cells.Item(25, 7).Value = 1.017114548628683
$ws.Cells.Item(25, 8).Value = 0.9783484700196539
$ws.Cells.Item(25, 12).Value = 0.2417418304359273
$ws.Cells.Item(25, 14).Value = 1.276865254151907
